$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the missing Status value (0) for row 3 (andy / Andy, THU20)
$ws.Range("D3").Value = 0

# Update the current selection/view to D8 (also clears the scrolled
# topLeftCell="A4" since the view now naturally shows D8 without scrolling)
[void]$ws.Range("D8").Select()
